# Data change after updating parser.py
# Update correlation analysis figures on Sheet1 to reflect the new parser output.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 ("Something - COG Dataset 1"): # of snippets with warnings 3 -> 4
$ws.Range("C2").Value = 4

# Row 3 ("Human Judgement (readability rating) - COG Dataset 3"):
# Kendall's Tau and Spearman's Rho recalculated by the updated parser
$ws.Range("F3").Value = -0.03063026539342525
$ws.Range("G3").Value = -0.03818805867014237

# Row 4 ("Something - fMRI Dataset"): # of warnings 0 -> 3
$ws.Range("D4").Value = 3
